# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" summary row at the top of the "总计" sheet
#    (pushing the existing quarters down by one row).
# 2. Insert a brand-new "2022-Q4" worksheet (with its own fund-holdings
#    table) right after "总计", pushing 2022-Q3 / 2022-Q2 / 2022-Q1 /
#    2021-Q4 / 2021-Q3 down one tab position each.
#
# NOTE: formatting is propagated with the single-argument-destination form
# of Range.Copy(Destination) rather than Copy() + PasteSpecial(). The
# clipboard-based two-step version was observed to silently drop the
# copied value/format when a new sheet had just been inserted via
# Worksheets.Add() earlier in the same script - Copy(Destination) copies
# directly cell-to-cell and was reliable in every trial.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert row 2 for 2022-Q4
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Give the new A2 the same "index column" look (bold/centered/bordered)
# as the rest of column A by copying formats+value from the row below.
$summary.Range("A3").Copy($summary.Range("A2"))

# The row-insert also bleeds formatting into B2:D2 - strip it back to
# the plain (unstyled) look the other data rows use.
$summary.Range("B2:D2").ClearFormats()

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.19

# Re-number the index column (A2:A7 = 0..5) now that everything shifted.
for ($i = 2; $i -le 7; $i++) {
  $summary.Range("A$i").Value = $i - 2
}

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q4" worksheet, placed right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Pull over the header-row / index-column formatting from the 2022-Q3
# sheet so the new tab matches the house style (bold, centered, bordered).
$template = $wb.Worksheets.Item("2022-Q3")
$template.Range("B1:H1").Copy($q4.Range("B1:H1"))
$template.Range("A2:A5").Copy($q4.Range("A2:A5"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B (fund code) and D:G (figures kept as text, same as the other
# quarterly sheets) must stay text so leading/trailing zeros survive.
$q4.Range("B2:B5").NumberFormat = "@"
$q4.Range("D2:G5").NumberFormat = "@"

$rows = @(
  @(0, "001637", "嘉实量化精选股票",             "12.82", "92.31", "1.18", "0.1513", 2),
  @(1, "562900", "易方达中证现代农业主题ETF",     "0.68",  "98.41", "2.21", "0.0150", 9),
  @(2, "008778", "嘉实中证500指数增强A",          "0.60",  "93.52", "1.88", "0.0113", 3),
  @(3, "008779", "嘉实中证500指数增强C",          "0.40",  "93.52", "1.88", "0.0075", 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = 2 + $i
  $d = $rows[$i]
  $q4.Range("A$r").Value = $d[0]
  $q4.Range("B$r").Value = $d[1]
  $q4.Range("C$r").Value = $d[2]
  $q4.Range("D$r").Value = $d[3]
  $q4.Range("E$r").Value = $d[4]
  $q4.Range("F$r").Value = $d[5]
  $q4.Range("G$r").Value = $d[6]
  $q4.Range("H$r").Value = $d[7]
}
